# FIX: Se corrige puntaje y se añaden resultados
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for columns C and D
$ws.Range("C1").Value = "Tiempo"
$ws.Range("D1").Value = "Entidades construidas"

# Corrected "Puntaje" values (column B) + new "Tiempo" (C) and
# "Entidades construidas" (D) data, rows 2..12
$data = @(
    @(0.2035947712418301, 623, 5),
    @(0.4218954248366013, 475, 2),
    @(0.2957516339869282, 633, 4),
    @(0.2709150326797386, 829, 5),
    @(0.3516339869281046, 532, 3),
    @(0.1895424836601307, 580, 5),
    @(0.4842528735632184, 631, 2),
    @(0.3577777777777778, 522, 3),
    @(0.4983141762452107, 414, 1),
    @(0.3677777777777778, 551, 3),
    @(0.4304597701149425, 475, 2)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 2).Value = $entry[0]
    $ws.Cells.Item($row, 3).Value = $entry[1]
    $ws.Cells.Item($row, 4).Value = $entry[2]
    $row++
}
